# Add "PAD number of interlace and double ring" comparison table
# on sheet "T28工艺库单元特性-面积" (3rd worksheet).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(3)

# ---------------------------------------------------------------
# Column width for column H
# ---------------------------------------------------------------
$ws.Range("H1").ColumnWidth = 12.21875

# ---------------------------------------------------------------
# Header row (row 10) - plan names, each merged across 3 columns
# ---------------------------------------------------------------
$ws.Range("I10:K10").Merge()
$ws.Range("I10").Value = "Plan A: 单圈PAD"

$ws.Range("L10:N10").Merge()
$ws.Range("L10").Value = "Plan B: 交错圈PAD"

$ws.Range("O10:Q10").Merge()
$ws.Range("O10").Value = "Plan B: 双圈PAD"

$ws.Range("R10:T10").Merge()
$ws.Range("R10").Value = "Plan C: 双圈交错PAD"

# ---------------------------------------------------------------
# Sub header row (row 11) - 长/宽/PAD个数 repeated per plan
# ---------------------------------------------------------------
$subHeaders = "长", "宽", "PAD个数"
$groups = "I", "L", "O", "R"
foreach ($g in $groups) {
    $cols = $ws.Range($g + "11").Resize(1, 3)
    for ($i = 0; $i -lt 3; $i++) {
        $cols.Cells.Item(1, $i + 1).Value = $subHeaders[$i]
    }
}

# ---------------------------------------------------------------
# Row labels (column H)
# ---------------------------------------------------------------
$ws.Range("H12").Value = "Chip Size"
$ws.Range("H13").Value = "PAD"
$ws.Range("H14").Value = "Core Size"
$ws.Range("H15").Value = "Core利用率"
$ws.Range("H16").Value = "PAD/Core"

# ---------------------------------------------------------------
# Data - Plan A (I/J/K)
# ---------------------------------------------------------------
$ws.Range("I12").Value = 1.5
$ws.Range("J12").Value = 2.2

$ws.Range("I13").Value = 0.25
$ws.Range("K13").Value = 83

$ws.Range("I14").Formula = "=I12-I13*2"
$ws.Range("J14").Formula = "=J12-I13*2"

$ws.Range("I15").Formula = "=I14*J14/I12/J12"

$ws.Range("I16").Formula = "=K13/I14/J14"

# ---------------------------------------------------------------
# Data - Plan B: 交错圈PAD (L/M/N)
# ---------------------------------------------------------------
$ws.Range("L12").Value = 1.5
$ws.Range("M12").Value = 2.2

$ws.Range("L13").Value = 0.35
$ws.Range("N13").Formula = "=K13*1.5"

$ws.Range("L14").Formula = "=L12-L13*2"
$ws.Range("M14").Formula = "=M12-L13*2"

$ws.Range("L15").Formula = "=L14*M14/L12/M12"

$ws.Range("L16").Formula = "=N13/L14/M14"

# ---------------------------------------------------------------
# Data - Plan B: 双圈PAD (O/P/Q)
# ---------------------------------------------------------------
$ws.Range("O12").Value = 1.5
$ws.Range("P12").Value = 2.2

$ws.Range("O13").Value = 0.5
$ws.Range("Q13").Formula = "=K13*1.8"

$ws.Range("O14").Formula = "=O12-O13*2"
$ws.Range("P14").Formula = "=P12-O13*2"

$ws.Range("O15").Formula = "=O14*P14/O12/P12"

$ws.Range("O16").Formula = "=Q13/O14/P14"

# ---------------------------------------------------------------
# Data - Plan C: 双圈交错PAD (R/S/T)
# ---------------------------------------------------------------
$ws.Range("R12").Value = 1.5
$ws.Range("S12").Value = 2.2

$ws.Range("R13").Value = 0.7
$ws.Range("T13").Formula = "=K13*1.5*1.8"

$ws.Range("R14").Formula = "=R12-R13*2"
$ws.Range("S14").Formula = "=S12-R13*2"

$ws.Range("R15").Formula = "=R14*S14/R12/S12"

$ws.Range("R16").Formula = "=T13/R14/S14"

# ---------------------------------------------------------------
# Highlight the "best" cell (Plan B 交错圈PAD utilization) in bold green,
# matching the style used elsewhere in this sheet to flag a stand-out value.
# ---------------------------------------------------------------
$ws.Range("L15").Font.Bold = $true
$ws.Range("L15").Font.Color = 5287936

# ---------------------------------------------------------------
# Alignment for data / header cells
# ---------------------------------------------------------------
$ws.Range("I10:T11").HorizontalAlignment = -4108
$ws.Range("I10:T11").VerticalAlignment = -4108
$ws.Range("I12:T16").HorizontalAlignment = -4108
$ws.Range("I12:T16").VerticalAlignment = -4108
$ws.Range("H12:H16").HorizontalAlignment = -4108
$ws.Range("H12:H16").VerticalAlignment = -4108

# ---------------------------------------------------------------
# Borders - outer box (medium) around H10:T15 per plan group, with
# thin separators inside each group, and a thin boxed row for row 16.
# ---------------------------------------------------------------

# Thin grid across the whole block first
$block = $ws.Range("H10:T16")
$block.Borders.Item(11).LineStyle = 1
$block.Borders.Item(11).Weight = 2
$block.Borders.Item(12).LineStyle = 1
$block.Borders.Item(12).Weight = 2
$block.Borders.Item(7).LineStyle = 1
$block.Borders.Item(7).Weight = 2
$block.Borders.Item(8).LineStyle = 1
$block.Borders.Item(8).Weight = 2
$block.Borders.Item(9).LineStyle = 1
$block.Borders.Item(9).Weight = 2
$block.Borders.Item(10).LineStyle = 1
$block.Borders.Item(10).Weight = 2

# Medium outer box around the main comparison table H10:T15
$mainBox = $ws.Range("H10:T15")
$mainBox.BorderAround(1, -4138)

# Medium separators between plan groups and around column H
$sepCols = "H", "K", "N", "Q"
foreach ($col in $sepCols) {
    $ws.Range($col + "10:" + $col + "15").Borders.Item(10).LineStyle = 1
    $ws.Range($col + "10:" + $col + "15").Borders.Item(10).Weight = -4138
}

# Medium line under the header row (row 10)
$ws.Range("H10:T10").Borders.Item(9).LineStyle = 1
$ws.Range("H10:T10").Borders.Item(9).Weight = -4138

# Row 16 (PAD/Core) thin boxed cells, separate from the box above
$ws.Range("H16:T16").BorderAround(1, 2)
$ws.Range("H16:T16").Borders.Item(11).LineStyle = 1
$ws.Range("H16:T16").Borders.Item(11).Weight = 2
$ws.Range("H16:T16").Borders.Item(12).LineStyle = 1
$ws.Range("H16:T16").Borders.Item(12).Weight = 2

# ---------------------------------------------------------------
# View settings - match author's final selection / scroll position
# ---------------------------------------------------------------
$ws.Range("K18").Select()
$excel.ActiveWindow.ScrollColumn = 2
